$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.257.54'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.926.07'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7170'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9991'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.86'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3192'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07090'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7911'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07970'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D13').Value = '1.924.97'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.389'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.83'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '30.234.97'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '256.74'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008071'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.764'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.84%  '
$ws.Range('D21').Value = '2.177.89'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.834'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.541'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.42'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.262'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1266'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.358'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.523'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.396'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.119'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05136'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.269'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7450'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.761'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01956'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.78%  '
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '77.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.365'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4508'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.990'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8444'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9984'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.62'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.722'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.427'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06125'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4199'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.43%  '
